# Apply the "preparation_temperature" -> "preparation_condition" and
# "storage_temperature" -> "storage_method" field renames, including the
# list sheets, their content, the data-validation formulas/messages, and
# the explanatory cell comments.

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Rename the two lookup-list sheets.
# ------------------------------------------------------------------
$wsPrepCond = $wb.Worksheets.Item("preparation_temperature list")
$wsPrepCond.Name = "preparation_condition list"

$wsStorageMethod = $wb.Worksheets.Item("storage_temperature list")
$wsStorageMethod.Name = "storage_method list"

# ------------------------------------------------------------------
# 2. Replace the "preparation_condition list" contents (was 8 rows,
#    now 7 rows).
# ------------------------------------------------------------------
$wsPrepCond.Range("A1").Value = "frozen in liquid nitrogen"
$wsPrepCond.Range("A2").Value = "frozen in liquid nitrogen vapor"
$wsPrepCond.Range("A3").Value = "frozen in ice"
$wsPrepCond.Range("A4").Value = "frozen in dry ice"
$wsPrepCond.Range("A5").Value = "frozen at -20 C"
$wsPrepCond.Range("A6").Value = "ambient temperature"
$wsPrepCond.Range("A7").Value = "unknown"
$wsPrepCond.Range("A8").ClearContents()

# ------------------------------------------------------------------
# 3. Replace the "storage_method list" contents (was 12 rows, now 11
#    rows).
# ------------------------------------------------------------------
$wsStorageMethod.Range("A1").Value = "frozen in liquid nitrogen"
$wsStorageMethod.Range("A2").Value = "frozen in liquid nitrogen vapor"
$wsStorageMethod.Range("A3").Value = "frozen in ice"
$wsStorageMethod.Range("A4").Value = "frozen in dry ice"
$wsStorageMethod.Range("A5").Value = "frozen at -80 C"
$wsStorageMethod.Range("A6").Value = "frozen at -20 C"
$wsStorageMethod.Range("A7").Value = "refrigerator"
$wsStorageMethod.Range("A8").Value = "ambient temperature"
$wsStorageMethod.Range("A9").Value = "incubated at 37 C"
$wsStorageMethod.Range("A10").Value = "none"
$wsStorageMethod.Range("A11").Value = "unknown"
$wsStorageMethod.Range("A12").ClearContents()

# ------------------------------------------------------------------
# 4. Update the header row field names and the data validation on the
#    main "Export as TSV" sheet for columns G (preparation_condition)
#    and K (storage_method).
# ------------------------------------------------------------------
$wsMain = $wb.Worksheets.Item("Export as TSV")

$wsMain.Range("G1").Value = "preparation_condition"
$wsMain.Range("K1").Value = "storage_method"

$dvG = $wsMain.Range("G2:G1048576").Validation
$dvG.Formula1 = "'preparation_condition list'!`$A`$1:`$A`$7"
$dvG.ErrorMessage = "Value must come from preparation_condition list."

$dvK = $wsMain.Range("K2:K1048576").Validation
$dvK.Formula1 = "'storage_method list'!`$A`$1:`$A`$11"
$dvK.ErrorMessage = "Value must come from storage_method list."

# ------------------------------------------------------------------
# 5. Update the header-row comments describing columns G and K.
# ------------------------------------------------------------------
$wsMain.Range("G1").Comment.Text("The condition under which the preparation occurred, such as whether the sample was placed in dry ice during the preparation.")
$wsMain.Range("K1").Comment.Text("The method by which the sample was stored, after preparation and before the assay was performed.")
